# Hackimon - Our team member
# Slide 1, "TextBox 19" ("Why we join?") bullet list:
#   1. Fix typo: "Hackidemy" -> "Hackademy"
#   2. Split the run " is a good learning program. " into " " + "is a good learning program. "
#   3. Merge the runs " to change the " + "world" back into a single run " to change the world"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$tb = $s.Shapes.Item("TextBox 19")
$tr = $tb.TextFrame.TextRange

# 1) Hackidemy -> Hackademy (simple in-place text fix, run stays intact)
$typo = $tr.Find("Hackidemy", 0)
if ($typo -ne $null) {
    $typo.Text = "Hackademy"
}

# 2) " is a good learning program. " -> split into two runs: " " and "is a good learning program. "
#    Re-writing just the leading space (a strict subset of the run) forces the
#    underlying run to split into two runs with identical formatting.
$sentence = $tr.Find("is a good learning program.", 0)
if ($sentence -ne $null) {
    $leadSpace = $tr.Characters($sentence.Start - 1, 1)
    $leadSpace.Text = " "
}

# 3) " to change the " + "world" -> merge back into a single run " to change the world"
#    Re-writing a range that spans both existing runs collapses them into one run.
$phrase = $tr.Find(" to change the world", 0)
if ($phrase -ne $null) {
    $phrase.Text = " to change the world"
}
